$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Move the two summary rows down: old row10 -> row11, old row11 -> row12
$ws.Range("D12").Formula = "=3*SUM(D2:D8)"
$ws.Range("C12").Value = "time to build"
$ws.Range("D11").Formula = "=D2-SUM(D3:D8)"
$ws.Range("E11").Formula = "=SUM(E2:E8)"
$ws.Range("C11").ClearContents()
$ws.Range("D10:E10").ClearContents()

# "move" count increases
$ws.Range("D2").Value = 8

# Row 5 becomes the new "attack" body part
$ws.Range("B5").Value = "attack"
$ws.Range("C5").Value = 80
$ws.Range("D5").Value = 3

# Row 6 becomes "rangedAttack"
$ws.Range("B6").Value = "rangedAttack"
$ws.Range("C6").Value = 150
$ws.Range("D6").Value = 0

# Row 7 becomes "heal"
$ws.Range("B7").Value = "heal"
$ws.Range("C7").Value = 250
$ws.Range("D7").Value = 0

# Row 8 becomes "tough" (new row)
$ws.Range("B8").Value = "tough"
$ws.Range("C8").Value = 10
$ws.Range("D8").Value = 5

# Re-fill the cost formula down the block, cell by cell so each one binds
# its own correct, independently-evaluable formula.
$ws.Range("E3").Formula = "=C3*D3"
$ws.Range("E4").Formula = "=C4*D4"
$ws.Range("E5").Formula = "=C5*D5"
$ws.Range("E6").Formula = "=C6*D6"
$ws.Range("E7").Formula = "=C7*D7"
$ws.Range("E8").Formula = "=C8*D8"

$ws.Range("D9").Select()
